$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.1350629056282652, 17.3793037085187, 0.4960521812499193),
    @(0.1016206307172993, 19.75333766694175, 0.4789792204747245),
    @(0.06108207437865381, 23.85121192697551, 0.2817078893979832),
    @(0.05634322691776208, 26.2486084269945, 0.494516332137253),
    @(0.06969296911682632, 26.07720179501735, 0.6387577184106538),
    @(0.06067078274357841, 29.85133070873766, 0.4812835823967427),
    @(0.1008969320459298, 31.6273879798011, 0.3881698754362048),
    @(0.103575790476888, 34.6009746266031, 0.6264967781249315),
    @(0.1262253520337998, 33.61862432298814, 0.6013241853193693),
    @(0.1218480528061232, 34.00647442646644, 0.7177278670250145),
    @(0.1046893883530775, 37.06859728901351, 0.4602523857549154),
    @(0.07514080095025165, 38.16035197397574, 0.663602089619127),
    @(0.1166433897191031, 39.45385403513613, 0.3882435901297552),
    @(0.09926335682753863, 37.85426367235398, 0.3078391237458981),
    @(0.05232827427921165, 38.65343111790397, 0.5455796909816264),
    @(0.102180339106266, 38.70765324435842, 0.610214476017145),
    @(0.1192574534390154, 38.35797666710676, 0.3101769623762196),
    @(0.1476867174960188, 37.01633085990378, 0.3300551408763464),
    @(0.07154082584924241, 35.99943341642707, 0.6545905976114401),
    @(0.06377467115624964, 35.55975371648374, 0.2916996486562654)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $data[$i][0]
    $ws.Range("C$row").Value = $data[$i][1]
    $ws.Range("D$row").Value = $data[$i][2]
}
